$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 6967
$ws.Cells.Item(80, 9).Value = 10200.5
$ws.Cells.Item(80, 10).Value = 500
$ws.Cells.Item(80, 11).Value = 30601.5
$ws.Cells.Item(80, 12).Value = 1500
$ws.Cells.Item(80, 13).Value = -29603.5
$ws.Cells.Item(80, 14).Value = -3496
$ws.Cells.Item(83, 8).Value = 6967
$ws.Cells.Item(83, 9).Value = 10200.5
$ws.Cells.Item(83, 10).Value = 500
$ws.Cells.Item(83, 11).Value = 91804.5
$ws.Cells.Item(83, 12).Value = 4500
$ws.Cells.Item(83, 13).Value = -86812.5
$ws.Cells.Item(83, 14).Value = -14484
$ws.Cells.Item(86, 8).Value = 46171150
$ws.Cells.Item(86, 9).Value = 120041220
$ws.Cells.Item(86, 10).Value = 2362.5
$ws.Cells.Item(86, 11).Value = 120041220
$ws.Cells.Item(86, 12).Value = 2362.5
$ws.Cells.Item(86, 13).Value = -120040097
$ws.Cells.Item(86, 14).Value = -4608.5
$ws.Cells.Item(89, 8).Value = 46171150
$ws.Cells.Item(89, 9).Value = 120041220
$ws.Cells.Item(89, 10).Value = 2362.5
$ws.Cells.Item(89, 11).Value = 600206100
$ws.Cells.Item(89, 12).Value = 11812.5
$ws.Cells.Item(89, 13).Value = -600200484
$ws.Cells.Item(89, 14).Value = -23044.5
$ws.Cells.Item(95, 8).Value = 50400
$ws.Cells.Item(95, 10).Value = 50400
$ws.Cells.Item(95, 12).Value = 50400
$ws.Cells.Item(95, 14).Value = -55892
$ws.Cells.Item(98, 8).Value = 680.2857
$ws.Cells.Item(98, 9).Value = 293.66666
$ws.Cells.Item(98, 11).Value = 293.66666
$ws.Cells.Item(98, 13).Value = 1204.33334
$ws.Cells.Item(115, 8).Value = 2165.4443
$ws.Cells.Item(115, 9).Value = 2047.5
$ws.Cells.Item(115, 10).Value = 2401.3333
$ws.Cells.Item(115, 11).Value = 6142.5
$ws.Cells.Item(115, 12).Value = 7203.999899999999
$ws.Cells.Item(115, 13).Value = -4575.5
$ws.Cells.Item(115, 14).Value = -10337.9999
$ws.Cells.Item(116, 8).Value = 2927.9167
$ws.Cells.Item(116, 9).Value = 2736.818
$ws.Cells.Item(116, 10).Value = 3089.6155
$ws.Cells.Item(116, 11).Value = 2736.818
$ws.Cells.Item(116, 12).Value = 3089.6155
$ws.Cells.Item(116, 13).Value = 705.1819999999998
$ws.Cells.Item(116, 14).Value = -9973.6155
$ws.Cells.Item(122, 8).Value = 680.2857
$ws.Cells.Item(122, 9).Value = 293.66666
$ws.Cells.Item(122, 11).Value = 880.9999799999999
$ws.Cells.Item(122, 13).Value = 1569.00002
$ws.Cells.Item(123, 8).Value = 28800
$ws.Cells.Item(123, 10).Value = 28800
$ws.Cells.Item(123, 12).Value = 28800
$ws.Cells.Item(123, 14).Value = -38600
$ws.Cells.Item(129, 8).Value = 996
$ws.Cells.Item(129, 9).Value = 504
$ws.Cells.Item(129, 10).Value = 1126.5306
$ws.Cells.Item(129, 11).Value = 1512
$ws.Cells.Item(129, 12).Value = 3379.5918
$ws.Cells.Item(129, 13).Value = 3488
$ws.Cells.Item(129, 14).Value = -13379.5918
$ws.Cells.Item(137, 8).Value = 3894.2917
$ws.Cells.Item(137, 9).Value = 3666.0625
$ws.Cells.Item(137, 10).Value = 4350.75
$ws.Cells.Item(137, 11).Value = 10998.1875
$ws.Cells.Item(137, 12).Value = 13052.25
$ws.Cells.Item(137, 13).Value = -8448.1875
$ws.Cells.Item(137, 14).Value = -18152.25
$ws.Cells.Item(138, 8).Value = 2164.532
$ws.Cells.Item(138, 9).Value = 1673.8148
$ws.Cells.Item(138, 10).Value = 2827
$ws.Cells.Item(138, 11).Value = 5021.4444
$ws.Cells.Item(138, 12).Value = 8481
$ws.Cells.Item(138, 13).Value = 118.5555999999997
$ws.Cells.Item(138, 14).Value = -18761

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 48907.855
$ws.Cells.Item(24, 10).Value = 48907.855
$ws.Cells.Item(24, 12).Value = 48907.855
$ws.Cells.Item(24, 14).Value = -49655.855
$ws.Cells.Item(69, 8).Value = 70000
$ws.Cells.Item(69, 10).Value = 70000
$ws.Cells.Item(69, 12).Value = 70000
$ws.Cells.Item(69, 14).Value = -71498
$ws.Cells.Item(72, 8).Value = 70000
$ws.Cells.Item(72, 10).Value = 70000
$ws.Cells.Item(72, 12).Value = 210000
$ws.Cells.Item(72, 14).Value = -217488
$ws.Cells.Item(74, 8).Value = 967.58826
$ws.Cells.Item(74, 9).Value = 744.3570999999999
$ws.Cells.Item(74, 10).Value = 2009.3334
$ws.Cells.Item(74, 11).Value = 744.3570999999999
$ws.Cells.Item(74, 12).Value = 2009.3334
$ws.Cells.Item(74, 13).Value = 129.6429000000001
$ws.Cells.Item(74, 14).Value = -3757.3334
$ws.Cells.Item(77, 8).Value = 967.58826
$ws.Cells.Item(77, 9).Value = 744.3570999999999
$ws.Cells.Item(77, 10).Value = 2009.3334
$ws.Cells.Item(77, 11).Value = 3721.7855
$ws.Cells.Item(77, 12).Value = 10046.667
$ws.Cells.Item(77, 13).Value = 646.2145
$ws.Cells.Item(77, 14).Value = -18782.667
$ws.Cells.Item(88, 8).Value = 2552.111
$ws.Cells.Item(88, 9).Value = 2156.3333
$ws.Cells.Item(88, 11).Value = 2156.3333
$ws.Cells.Item(88, 13).Value = -1750.3333
$ws.Cells.Item(91, 8).Value = 2552.111
$ws.Cells.Item(91, 9).Value = 2156.3333
$ws.Cells.Item(91, 11).Value = 2156.3333
$ws.Cells.Item(91, 13).Value = -752.3332999999998
$ws.Cells.Item(100, 8).Value = 48907.855
$ws.Cells.Item(100, 10).Value = 48907.855
$ws.Cells.Item(100, 12).Value = 48907.855
$ws.Cells.Item(100, 14).Value = -51071.855
$ws.Cells.Item(102, 8).Value = 4208.0586
$ws.Cells.Item(102, 9).Value = 4346.0625
$ws.Cells.Item(102, 11).Value = 4346.0625
$ws.Cells.Item(102, 13).Value = -2724.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 3463
$ws.Cells.Item(22, 9).Value = 3463
$ws.Cells.Item(22, 11).Value = 3463
$ws.Cells.Item(22, 13).Value = -3290
$ws.Cells.Item(107, 8).Value = 201596
$ws.Cells.Item(107, 9).Value = 334460
$ws.Cells.Item(107, 11).Value = 334460
$ws.Cells.Item(107, 13).Value = -332540

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4551.757
$ws.Cells.Item(31, 9).Value = 1102.5853
$ws.Cells.Item(31, 10).Value = 8837.091
$ws.Cells.Item(31, 11).Value = 1102.5853
$ws.Cells.Item(31, 12).Value = 8837.091
$ws.Cells.Item(31, 13).Value = -807.5853
$ws.Cells.Item(31, 14).Value = -9427.091
$ws.Cells.Item(34, 8).Value = 4551.757
$ws.Cells.Item(34, 9).Value = 1102.5853
$ws.Cells.Item(34, 10).Value = 8837.091
$ws.Cells.Item(34, 11).Value = 1102.5853
$ws.Cells.Item(34, 12).Value = 8837.091
$ws.Cells.Item(34, 13).Value = -900.5853
$ws.Cells.Item(34, 14).Value = -9241.091
$ws.Cells.Item(107, 8).Value = 1076.625
$ws.Cells.Item(107, 9).Value = 500
$ws.Cells.Item(107, 11).Value = 500
$ws.Cells.Item(107, 13).Value = 1420
$ws.Cells.Item(134, 8).Value = 7590.1577
$ws.Cells.Item(134, 9).Value = 8922.615
$ws.Cells.Item(134, 10).Value = 4703.1665
$ws.Cells.Item(134, 11).Value = 26767.845
$ws.Cells.Item(134, 12).Value = 14109.4995
$ws.Cells.Item(134, 13).Value = -24232.845
$ws.Cells.Item(134, 14).Value = -19179.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 166.66667
$ws.Cells.Item(2, 9).Value = 90
$ws.Cells.Item(2, 10).Value = 205
$ws.Cells.Item(2, 11).Value = 540
$ws.Cells.Item(2, 12).Value = 1230
$ws.Cells.Item(2, 13).Value = -427
$ws.Cells.Item(2, 14).Value = -1456
$ws.Cells.Item(17, 8).Value = 920.1
$ws.Cells.Item(17, 9).Value = 533.3333
$ws.Cells.Item(17, 10).Value = 1500.25
$ws.Cells.Item(17, 11).Value = 1599.9999
$ws.Cells.Item(17, 12).Value = 4500.75
$ws.Cells.Item(17, 13).Value = -1430.9999
$ws.Cells.Item(17, 14).Value = -4838.75
$ws.Cells.Item(62, 8).Value = 9006.833000000001
$ws.Cells.Item(62, 10).Value = 9006.833000000001
$ws.Cells.Item(62, 12).Value = 27020.499
$ws.Cells.Item(62, 14).Value = -28392.499
$ws.Cells.Item(65, 8).Value = 9006.833000000001
$ws.Cells.Item(65, 10).Value = 9006.833000000001
$ws.Cells.Item(65, 12).Value = 81061.497
$ws.Cells.Item(65, 14).Value = -87925.497

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 19966.334
$ws.Cells.Item(39, 10).Value = 19966.334
$ws.Cells.Item(39, 12).Value = 19966.334
$ws.Cells.Item(39, 14).Value = -21030.334
$ws.Cells.Item(80, 8).Value = 36387590
$ws.Cells.Item(80, 9).Value = 84834750
$ws.Cells.Item(80, 10).Value = 52221
$ws.Cells.Item(80, 11).Value = 84834750
$ws.Cells.Item(80, 12).Value = 52221
$ws.Cells.Item(80, 13).Value = -84833752
$ws.Cells.Item(80, 14).Value = -54217
$ws.Cells.Item(83, 8).Value = 36387590
$ws.Cells.Item(83, 9).Value = 84834750
$ws.Cells.Item(83, 10).Value = 52221
$ws.Cells.Item(83, 11).Value = 424173750
$ws.Cells.Item(83, 12).Value = 261105
$ws.Cells.Item(83, 13).Value = -424168758
$ws.Cells.Item(83, 14).Value = -271089
$ws.Cells.Item(95, 8).Value = 60248.8
$ws.Cells.Item(95, 10).Value = 60248.8
$ws.Cells.Item(95, 12).Value = 60248.8
$ws.Cells.Item(95, 14).Value = -65740.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 8393.823
$ws.Cells.Item(22, 9).Value = 1880
$ws.Cells.Item(22, 10).Value = 11107.917
$ws.Cells.Item(22, 11).Value = 1880
$ws.Cells.Item(22, 12).Value = 11107.917
$ws.Cells.Item(22, 13).Value = -1585
$ws.Cells.Item(22, 14).Value = -11697.917
$ws.Cells.Item(27, 8).Value = 8393.823
$ws.Cells.Item(27, 9).Value = 1880
$ws.Cells.Item(27, 10).Value = 11107.917
$ws.Cells.Item(27, 11).Value = 1880
$ws.Cells.Item(27, 12).Value = 11107.917
$ws.Cells.Item(27, 13).Value = -1773
$ws.Cells.Item(27, 14).Value = -11321.917
$ws.Cells.Item(82, 8).Value = 50002456
$ws.Cells.Item(82, 9).Value = 71431370
$ws.Cells.Item(82, 10).Value = 1656.3334
$ws.Cells.Item(82, 11).Value = 71431370
$ws.Cells.Item(82, 12).Value = 1656.3334
$ws.Cells.Item(82, 13).Value = -71431009
$ws.Cells.Item(82, 14).Value = -2378.3334
$ws.Cells.Item(85, 8).Value = 50002456
$ws.Cells.Item(85, 9).Value = 71431370
$ws.Cells.Item(85, 10).Value = 1656.3334
$ws.Cells.Item(85, 11).Value = 71431370
$ws.Cells.Item(85, 12).Value = 1656.3334
$ws.Cells.Item(85, 13).Value = -71430122
$ws.Cells.Item(85, 14).Value = -4152.3334
$ws.Cells.Item(122, 8).Value = 3648.25
$ws.Cells.Item(122, 9).Value = 3363.6365
$ws.Cells.Item(122, 10).Value = 3996.111
$ws.Cells.Item(122, 11).Value = 10090.9095
$ws.Cells.Item(122, 12).Value = 11988.333
$ws.Cells.Item(122, 13).Value = -7640.9095
$ws.Cells.Item(122, 14).Value = -16888.333
$ws.Cells.Item(132, 8).Value = 3215.5264
$ws.Cells.Item(132, 9).Value = 2269.125
$ws.Cells.Item(132, 10).Value = 3903.818
$ws.Cells.Item(132, 11).Value = 6807.375
$ws.Cells.Item(132, 12).Value = 11711.454
$ws.Cells.Item(132, 13).Value = -4277.375
$ws.Cells.Item(132, 14).Value = -16771.454
$ws.Cells.Item(136, 8).Value = 2499.889
$ws.Cells.Item(136, 9).Value = 2880.8
$ws.Cells.Item(136, 10).Value = 2023.75
$ws.Cells.Item(136, 11).Value = 8642.400000000001
$ws.Cells.Item(136, 12).Value = 6071.25
$ws.Cells.Item(136, 13).Value = -6092.400000000001
$ws.Cells.Item(136, 14).Value = -11171.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 64514.25
$ws.Cells.Item(31, 10).Value = 64514.25
$ws.Cells.Item(31, 12).Value = 64514.25
$ws.Cells.Item(31, 14).Value = -65210.25
$ws.Cells.Item(132, 8).Value = 1209.4237
$ws.Cells.Item(132, 9).Value = 898.907
$ws.Cells.Item(132, 10).Value = 2043.9375
$ws.Cells.Item(132, 11).Value = 2696.721
$ws.Cells.Item(132, 12).Value = 6131.8125
$ws.Cells.Item(132, 13).Value = -166.721
$ws.Cells.Item(132, 14).Value = -11191.8125
$ws.Cells.Item(136, 8).Value = 1918.2273
$ws.Cells.Item(136, 9).Value = 1646.0769
$ws.Cells.Item(136, 10).Value = 2311.3333
$ws.Cells.Item(136, 11).Value = 2009.3334
$ws.Cells.Item(136, 12).Value = 6933.999899999999
$ws.Cells.Item(136, 13).Value = -2388.2307
$ws.Cells.Item(136, 14).Value = -12033.9999
